$wb = $excel.ActiveWorkbook

# --- Sheet 1: ALC ---
$ws = $wb.Worksheets.Item(1)
$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()

$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()

$ws.Range("H116").Value = 2997
$ws.Range("I116").Value = 2997
$ws.Range("K116").Value = 2997
$ws.Range("M116").Value = 445

$ws.Range("H138").Value = 2514.2856

# --- Sheet 2: ARM ---
$ws = $wb.Worksheets.Item(2)
$ws.Range("H132").Value = 4678.5
$ws.Range("I132").Value = 4414.2
$ws.Range("J132").Value = 6000
$ws.Range("K132").Value = 13242.6
$ws.Range("L132").Value = 18000
$ws.Range("M132").Value = -10712.6
$ws.Range("N132").Value = -23060

# --- Sheet 3: BSM ---
$ws = $wb.Worksheets.Item(3)
$ws.Range("H38").Value = 7000
$ws.Range("J38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("N38").ClearContents()

$ws.Range("H86").Value = 0
$ws.Range("I86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("M86").ClearContents()

$ws.Range("H89").Value = 0
$ws.Range("I89").Value = 0
$ws.Range("K89").Value = 10000
$ws.Range("M89").ClearContents()

$ws.Range("H105").Value = 1000
$ws.Range("I105").Value = 1000
$ws.Range("K105").Value = 1000
$ws.Range("M105").Value = 747

# --- Sheet 4: CRP ---
$ws = $wb.Worksheets.Item(4)
$ws.Range("H4").Value = 601
$ws.Range("I4").Value = 1000
$ws.Range("J4").Value = 202
$ws.Range("K4").Value = 1000
$ws.Range("L4").Value = 202
$ws.Range("M4").Value = -888
$ws.Range("N4").Value = -426

$ws.Range("H7").Value = 39.3
$ws.Range("I7").Value = 11.4
$ws.Range("J7").Value = 67.2
$ws.Range("K7").Value = 11.4
$ws.Range("L7").Value = 67.2
$ws.Range("M7").Value = 101.6
$ws.Range("N7").Value = -293.2

$ws.Range("H22").Value = 949.5
$ws.Range("J22").Value = 1000
$ws.Range("L22").Value = 1000
$ws.Range("N22").Value = -1700

$ws.Range("H107").Value = 5000
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 5000
$ws.Range("K107").Value = 0
$ws.Range("L107").Value = 5000
$ws.Range("M107").ClearContents()
$ws.Range("N107").Value = -8840

$ws.Range("H124").Value = 148999.5
$ws.Range("J124").Value = 148999.5
$ws.Range("L124").Value = 148999.5
$ws.Range("N124").Value = -153909.5

$ws.Range("H141").Value = 1234258
$ws.Range("J141").Value = 1234258
$ws.Range("L141").Value = 1234258
$ws.Range("N141").Value = -1244618

# --- Sheet 5: CUL ---
$ws = $wb.Worksheets.Item(5)
$ws.Range("H8").Value = 1401.5
$ws.Range("I8").Value = 1401.5
$ws.Range("K8").Value = 4204.5
$ws.Range("M8").Value = -4065.5

$ws.Range("H40").Value = 46
$ws.Range("J40").Value = 49
$ws.Range("L40").Value = 196
$ws.Range("N40").Value = -334

$ws.Range("H46").Value = 186.25
$ws.Range("I46").Value = 95
$ws.Range("K46").Value = 285
$ws.Range("M46").Value = -194

$ws.Range("H108").Value = 565.8333
$ws.Range("I108").Value = 565.8333
$ws.Range("K108").Value = 1697.4999
$ws.Range("M108").Value = 1182.5001

$ws.Range("H131").Value = 2314.9167
$ws.Range("I131").Value = 1518.4286
$ws.Range("J131").Value = 3430
$ws.Range("K131").Value = 4555.2858
$ws.Range("L131").Value = 10290
$ws.Range("M131").Value = 484.7142000000003
$ws.Range("N131").Value = -20370

$ws.Range("H137").Value = 4008.25
$ws.Range("J137").Value = 4008.25
$ws.Range("L137").Value = 12024.75
$ws.Range("N137").Value = -22224.75

# --- Sheet 6: GSM ---
$ws = $wb.Worksheets.Item(6)
$ws.Range("H31").Value = 1396
$ws.Range("I31").Value = 1396
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 1396
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -1104
$ws.Range("N31").ClearContents()

$ws.Range("H37").Value = 1396
$ws.Range("I37").Value = 1396
$ws.Range("J37").Value = 0
$ws.Range("K37").Value = 1396
$ws.Range("L37").Value = 0
$ws.Range("M37").Value = -1119
$ws.Range("N37").ClearContents()

$ws.Range("H70").Value = 10500
$ws.Range("J70").Value = 12000
$ws.Range("L70").Value = 12000
$ws.Range("N70").Value = -12540

$ws.Range("H73").Value = 10500
$ws.Range("J73").Value = 12000
$ws.Range("L73").Value = 12000
$ws.Range("N73").Value = -13872

$ws.Range("H122").Value = 582.61536
$ws.Range("I122").Value = 582.61536
$ws.Range("K122").Value = 1747.84608
$ws.Range("M122").Value = 702.15392

$ws.Range("H126").Value = 1994
$ws.Range("J126").Value = 1994
$ws.Range("L126").Value = 5982
$ws.Range("N126").Value = -10922

$ws.Range("H132").Value = 4084.6667
$ws.Range("I132").Value = 3836.3333
$ws.Range("J132").Value = 4333
$ws.Range("K132").Value = 11508.9999
$ws.Range("L132").Value = 12999
$ws.Range("M132").Value = -8978.999899999999
$ws.Range("N132").Value = -18059

# --- Sheet 7: LTW ---
$ws = $wb.Worksheets.Item(7)
$ws.Range("H22").Value = 666
$ws.Range("I22").Value = 499.5
$ws.Range("J22").Value = 999
$ws.Range("K22").Value = 499.5
$ws.Range("L22").Value = 999
$ws.Range("M22").Value = -204.5
$ws.Range("N22").Value = -1589

$ws.Range("H27").Value = 666
$ws.Range("I27").Value = 499.5
$ws.Range("J27").Value = 999
$ws.Range("K27").Value = 499.5
$ws.Range("L27").Value = 999
$ws.Range("M27").Value = -392.5
$ws.Range("N27").Value = -1213

$ws.Range("H46").Value = 0
$ws.Range("I46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("M46").ClearContents()

$ws.Range("H55").Value = 596.75
$ws.Range("I55").Value = 534.61536
$ws.Range("J55").Value = 866
$ws.Range("K55").Value = 534.61536
$ws.Range("L55").Value = 866
$ws.Range("M55").Value = -361.61536
$ws.Range("N55").Value = -1212

$ws.Range("H132").Value = 5813.875
$ws.Range("I132").Value = 5101.8
$ws.Range("K132").Value = 15305.4
$ws.Range("M132").Value = -12775.4

$ws.Range("H136").Value = 57192
$ws.Range("I136").Value = 23992.182
$ws.Range("J136").Value = 97769.55499999999
$ws.Range("K136").Value = 71976.546
$ws.Range("L136").Value = 293308.665
$ws.Range("M136").Value = -69426.546
$ws.Range("N136").Value = -298408.665

# --- Sheet 8: WVR ---
$ws = $wb.Worksheets.Item(8)
$ws.Range("H49").Value = 0
$ws.Range("I49").Value = 0
$ws.Range("K49").Value = 0
$ws.Range("M49").ClearContents()

$ws.Range("H132").Value = 1220.2727
$ws.Range("I132").Value = 1142.3
$ws.Range("J132").Value = 2000
$ws.Range("K132").Value = 3426.9
$ws.Range("L132").Value = 6000
$ws.Range("M132").Value = -896.8999999999996
$ws.Range("N132").Value = -11060

$ws.Range("H136").Value = 2070
$ws.Range("I136").Value = 2149.0908
$ws.Range("K136").Value = 6447.2724
$ws.Range("M136").Value = -3897.2724
